$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.665.18'
$ws.Range('E2').Value = '  -0.22%  '
$ws.Range('D3').Value = '2.050.64'
$ws.Range('E3').Value = '  -0.40%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '246.52'
$ws.Range('E5').Value = '  +0.75%  '
$ws.Range('E6').Value = '  +2.68%  '
$ws.Range('D7').Value = '57.25'
$ws.Range('E7').Value = '  +3.19%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '63.36'
$ws.Range('E9').Value = '  +7.94%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.370'
$ws.Range('E10').Value = '  +1.30%  '
$ws.Range('D11').Value = '0.0755'
$ws.Range('E11').Value = '  -1.17%  '
$ws.Range('E12').Value = '  -2.63%  '
$ws.Range('D13').Value = '0.929'
$ws.Range('E13').Value = '  +7.06%  '
$ws.Range('D14').Value = '14.47'
$ws.Range('E14').Value = '  -2.56%  '
$ws.Range('D15').Value = '2.348.39'
$ws.Range('E15').Value = '  -0.17%  '
$ws.Range('D16').Value = '5.45'
$ws.Range('E16').Value = '  -1.17%  '
$ws.Range('D17').Value = '2.070.17'
$ws.Range('E17').Value = '  +0.74%  '
$ws.Range('D18').Value = '17.94'
$ws.Range('E18').Value = '  +3.92%  '
$ws.Range('D19').Value = '36.590.50'
$ws.Range('E19').Value = '  -0.35%  '
$ws.Range('D20').Value = '72.03'
$ws.Range('E20').Value = '  -1.00%  '
$ws.Range('D21').Value = '0.0₃0862'
$ws.Range('E21').Value = '  -0.99%  '
$ws.Range('D22').Value = '237.59'
$ws.Range('E22').Value = '  +1.08%  '
$ws.Range('D23').Value = '5.22'
$ws.Range('E23').Value = '  -2.89%  '
$ws.Range('E24').Value = '  +0.28%  '
$ws.Range('E25').Value = '  -2.03%  '
$ws.Range('D26').Value = '2.28'
$ws.Range('E26').Value = '  +4.17%  '
$ws.Range('D27').Value = '9.37'
$ws.Range('E27').Value = '  -5.10%  '
$ws.Range('D28').Value = '164.85'
$ws.Range('E28').Value = '  -0.99%  '
$ws.Range('D29').Value = '20.05'
$ws.Range('E29').Value = '  -1.55%  '
$ws.Range('E30').Value = '  -0.92%  '
$ws.Range('D31').Value = '1.21'
$ws.Range('E31').Value = '  +4.40%  '
$ws.Range('E32').Value = '  -5.41%  '
$ws.Range('E33').Value = '  -0.68%  '
$ws.Range('D34').Value = '4.46'
$ws.Range('E34').Value = '  -5.59%  '
$ws.Range('D35').Value = '0.0876'
$ws.Range('E35').Value = '  +3.55%  '
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('E37').Value = '  -0.76%  '
$ws.Range('D38').Value = '2.22'
$ws.Range('E38').Value = '  -4.77%  '
$ws.Range('E39').Value = '  +4.36%  '
$ws.Range('E40').Value = '  -4.38%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '0.0217'
$ws.Range('E41').Value = '  -0.94%  '
$ws.Range('B42').Value = 'HuobiToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.90'
$ws.Range('E42').Value = '  -0.93%  '
$ws.Range('E43').Value = '  -2.60%  '
$ws.Range('D44').Value = '94.36'
$ws.Range('E44').Value = '  -1.38%  '
$ws.Range('D45').Value = '0.0915'
$ws.Range('E45').Value = '  -3.67%  '
$ws.Range('D46').Value = '16.05'
$ws.Range('E46').Value = '  -0.96%  '
$ws.Range('D47').Value = '1.382.24'
$ws.Range('E47').Value = '  +5.00%  '
$ws.Range('D48').Value = '7.43'
$ws.Range('E48').Value = '  +7.76%  '
$ws.Range('D49').Value = '2.95'
$ws.Range('E49').Value = '  +3.07%  '
$ws.Range('D50').Value = '2.28'
$ws.Range('E50').Value = '  -2.35%  '
$ws.Range('D51').Value = '46.08'
$ws.Range('E51').Value = '  +1.79%  '
